# Applies the cryptos.xlsx data refresh described in the commit:
# "Updated cryptos list on Tue May 30 09:01:04 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain-text numbers (e.g. "27.821.12", "0.5030") that must
# stay text instead of being auto-coerced to floating point numbers, which
# would silently rewrite values like "0.9065" -> "0.90649999999999997" and
# drop significant trailing/leading zeros. Force the column to Text format
# before writing, then restore the default "Normal" style so no stray
# number-format style is left attached to the cells.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

$ws.Range("D2").Value = "27.821.12"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "1.905.06"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "313.22"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").Value = "0.5030"
$ws.Range("E7").Value = "  +4.55%  "
$ws.Range("D8").Value = "0.3808"
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").Value = "0.07269"
$ws.Range("E9").Value = "  -1.17%  "
$ws.Range("D10").Value = "0.9073"
$ws.Range("E10").Value = "  -2.60%  "
$ws.Range("D11").Value = "20.88"
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.07650"
$ws.Range("E12").Value = "  -1.28%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.897.59"
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("D14").Value = "5.492"
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("D15").Value = "91.72"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("D17").Value = "0.000008710"
$ws.Range("E17").Value = "  -1.28%  "
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").Value = "27.870.55"
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("D20").Value = "14.59"
$ws.Range("E20").Value = "  -1.31%  "
$ws.Range("D21").Value = "5.163"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").Value = "10.80"
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("D23").Value = "6.586"
$ws.Range("E23").Value = "  -0.56%  "
$ws.Range("D24").Value = "153.83"
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("D25").Value = "1.883"
$ws.Range("E25").Value = "  -1.93%  "
$ws.Range("D26").Value = "2.216"
$ws.Range("E26").Value = "  +4.17%  "
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("D28").Value = "115.40"
$ws.Range("E28").Value = "  -1.33%  "
$ws.Range("D29").Value = "4.908"
$ws.Range("E29").Value = "  -0.99%  "
$ws.Range("D30").Value = "0.09034"
$ws.Range("E30").Value = "  +1.09%  "
$ws.Range("D31").Value = "3.211"
$ws.Range("E31").Value = "  -2.55%  "
$ws.Range("D32").Value = "1.221"
$ws.Range("E32").Value = "  -2.99%  "
$ws.Range("D33").Value = "4.740"
$ws.Range("E33").Value = "  +1.39%  "
$ws.Range("D34").Value = "0.7693"
$ws.Range("E34").Value = "  -0.41%  "
$ws.Range("E35").Value = "  +0.53%  "
$ws.Range("D36").Value = "2.515"
$ws.Range("E36").Value = "  -4.45%  "
$ws.Range("E37").Value = "  -1.51%  "
$ws.Range("D38").Value = "0.5523"
$ws.Range("E38").Value = "  +0.72%  "
$ws.Range("D39").Value = "3.015"
$ws.Range("E39").Value = "  +0.81%  "
$ws.Range("E40").Value = "  -0.78%  "
$ws.Range("D41").Value = "6.887"
$ws.Range("E41").Value = "  -1.88%  "
$ws.Range("D42").Value = "8.465"
$ws.Range("E42").Value = "  -0.14%  "
$ws.Range("E43").Value = "  -1.02%  "
$ws.Range("D44").Value = "111.74"
$ws.Range("E44").Value = "  +3.64%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.4813"
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "10.54"
$ws.Range("E46").Value = "  -1.21%  "
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("D48").Value = "1.630"
$ws.Range("E48").Value = "  -1.14%  "
$ws.Range("D49").Value = "67.57"
$ws.Range("E49").Value = "  -0.52%  "
$ws.Range("D50").Value = "0.06068"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("D51").Value = "0.9039"
$ws.Range("E51").Value = "  +0.56%  "

$priceCol.Style = "Normal"
